$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The text that used to live in B2 ("MONTEREY AREA TOTALS") now lives in A2,
# and B2 becomes the much shorter "Totals" label (matching the pattern used
# by the other port subtotal rows, e.g. B14/B28/B36).
$ws.Range("A2").Value = "MONTEREY AREA TOTALS"
$ws.Range("B2").Value = "Totals"

# Column A's best-fit width now needs to be recalculated since it holds the
# long "MONTEREY AREA TOTALS" string (previously in column B), and column B's
# best-fit width is recalculated too (its longest remaining entry keeps the
# same effective width).
$ws.Columns("A:A").AutoFit()
$ws.Columns("B:B").AutoFit()

# Move the active selection to B6.
$ws.Range("B6").Select()
